# updated legacy GSC export data
#
# The "Chart" sheet (sheet1) holds a rolling window of daily export rows:
# column A = date label (text, not a real date), B/C = counts.
# This update drops the oldest day (2025-10-03, row 2) - which slides every
# other row up by one - and appends two new trailing days
# (2026-01-01, 2026-01-02) with zeroed counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest dated row; Excel shifts rows 3..91 up into 2..90,
# carrying each row's own B/C values along with it (matches the diff's
# per-row C-column shift).
$ws.Rows.Item(2).Delete()

# Helper: write a literal, date-formatted-looking string as TEXT (not an
# auto-converted date serial) without touching any cell's NumberFormat/style.
# Build the literal via a text formula, then PasteSpecial just the value
# (xlPasteValues = -4163) on top of itself so the stored cell becomes a
# plain shared-string literal, matching the source file's encoding.
function Set-TextValue($cell, [string]$text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Append the two new trailing rows at 91 and 92.
Set-TextValue $ws.Cells.Item(91, 1) "2026-01-01"
$ws.Cells.Item(91, 2).Value = 0.0
$ws.Cells.Item(91, 3).Value = 0.0

Set-TextValue $ws.Cells.Item(92, 1) "2026-01-02"
$ws.Cells.Item(92, 2).Value = 0.0
$ws.Cells.Item(92, 3).Value = 0.0

Write-Output "GSC export data updated"
